$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 489, shifting existing data down (old 489-504 -> 493-508)
$ws.Rows("489:492").Insert()

# Row 489
$ws.Cells.Item(489, 1).Value = 1
$ws.Cells.Item(489, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(489, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(489, 4).Value = 44568
$ws.Cells.Item(489, 5).Value = 15
$ws.Cells.Item(489, 6).Value = 100112024
$ws.Cells.Item(489, 7).Value = "Choclo"
$ws.Cells.Item(489, 8).Value = "Dulce o Americano"
$ws.Cells.Item(489, 9).Value = "Primera"
$ws.Cells.Item(489, 10).Value = 130
$ws.Cells.Item(489, 11).Value = 7000
$ws.Cells.Item(489, 12).Value = 8000
$ws.Cells.Item(489, 13).Value = 7500
$ws.Cells.Item(489, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(489, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(489, 16).Value = 107
$ws.Cells.Item(489, 17).Value = 70
$ws.Cells.Item(489, 18).Value = "Hortaliza"

# Row 490
$ws.Cells.Item(490, 1).Value = 1
$ws.Cells.Item(490, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(490, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(490, 4).Value = 44568
$ws.Cells.Item(490, 5).Value = 15
$ws.Cells.Item(490, 6).Value = 100112024
$ws.Cells.Item(490, 7).Value = "Choclo"
$ws.Cells.Item(490, 8).Value = "Lluteño"
$ws.Cells.Item(490, 9).Value = "Primera"
$ws.Cells.Item(490, 10).Value = 50
$ws.Cells.Item(490, 11).Value = 20000
$ws.Cells.Item(490, 12).Value = 22000
$ws.Cells.Item(490, 13).Value = 21000
$ws.Cells.Item(490, 14).Value = "$/saco 50 unidades"
$ws.Cells.Item(490, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(490, 16).Value = 420
$ws.Cells.Item(490, 17).Value = 50
$ws.Cells.Item(490, 18).Value = "Hortaliza"

# Row 491
$ws.Cells.Item(491, 1).Value = 1
$ws.Cells.Item(491, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(491, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(491, 4).Value = 44568
$ws.Cells.Item(491, 5).Value = 15
$ws.Cells.Item(491, 6).Value = 100112024
$ws.Cells.Item(491, 7).Value = "Choclo"
$ws.Cells.Item(491, 8).Value = "Lluteño"
$ws.Cells.Item(491, 9).Value = "Segunda"
$ws.Cells.Item(491, 10).Value = 50
$ws.Cells.Item(491, 11).Value = 17000
$ws.Cells.Item(491, 12).Value = 19000
$ws.Cells.Item(491, 13).Value = 18000
$ws.Cells.Item(491, 14).Value = "$/saco 75 unidades"
$ws.Cells.Item(491, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(491, 16).Value = 240
$ws.Cells.Item(491, 17).Value = 75
$ws.Cells.Item(491, 18).Value = "Hortaliza"

# Row 492
$ws.Cells.Item(492, 1).Value = 1
$ws.Cells.Item(492, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(492, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(492, 4).Value = 44568
$ws.Cells.Item(492, 5).Value = 15
$ws.Cells.Item(492, 6).Value = 100112024
$ws.Cells.Item(492, 7).Value = "Choclo"
$ws.Cells.Item(492, 8).Value = "Lluteño"
$ws.Cells.Item(492, 9).Value = "Tercera"
$ws.Cells.Item(492, 10).Value = 50
$ws.Cells.Item(492, 11).Value = 13000
$ws.Cells.Item(492, 12).Value = 15000
$ws.Cells.Item(492, 13).Value = 14000
$ws.Cells.Item(492, 14).Value = "$/saco 100 unidades"
$ws.Cells.Item(492, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(492, 16).Value = 140
$ws.Cells.Item(492, 17).Value = 100
$ws.Cells.Item(492, 18).Value = "Hortaliza"
